$d = $word.ActiveDocument

function ReplaceInParagraph($paraIndex, $oldText, $newText) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $ok = $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $ok) {
        Write-Output "MISS paragraph $paraIndex old=[$oldText]"
    }
}

# 1. "English" inside hyperlink (paragraph 1)
ReplaceInParagraph 1 "English" "ඉංග්‍රීසි"

# 2. "English" paragraph (P68B1DB1-Normal2 style) (paragraph 3)
ReplaceInParagraph 3 "English" "ඉංග්‍රීසි"

# 3. "Brief" (bold run) in table 1 (paragraph 5)
ReplaceInParagraph 5 "Brief" "සාරාංශය"

# 4. Brief description paragraph (paragraph 6)
ReplaceInParagraph 6 "An email sent to partners in the target country who have RSVPed yes. We want them to submit their documents. It will be sent via customer.io" "ඉතා කැඩපතක් සහිත ඊ-මේල් පණිවුඩයක් විශේෂිත රටේ පිළිගෙන ඇති ආරාධිත සහකාරයන්ට එවනය. ඔවුන්ගේ ලේඛන යැවීමට අපි කැමතියි. It will be sent via customer.io"

# 5. "Target audience" (bold run) (paragraph 8)
ReplaceInParagraph 8 "Target audience" "ඉලක්කය"

# 6. "Invited partners who RSVP yes" (paragraph 9)
ReplaceInParagraph 9 "Invited partners who RSVP yes" "ඉතා කැඩපතක් ලබාගත් ආරාධිත සහකාරයන්"

# 7. "Subject: " (bold run) -> "විෂය:" (bold) + new run " " (plain) (paragraph 12)
ReplaceInParagraph 12 "Subject: " "විෂය:"
$p12 = $d.Paragraphs(12)
$subjRange = $p12.Range.Duplicate
$subjRange.Find.Execute("විෂය:") | Out-Null
$subjRange.Collapse(0)
$subjRange.InsertAfter(" ")
$subjRange.Font.Bold = 0
# 8. " — take the next step" -> " — ඊළඟ පියවර ගන්න" (paragraph 12)
ReplaceInParagraph 12 " — take the next step" " — ඊළඟ පියවර ගන්න"

# 9. "Hi " (paragraph 15)
ReplaceInParagraph 15 "Hi " "ආයුබෝවන් "

# 10. "To confirm your registration..." (paragraph 18)
ReplaceInParagraph 18 "To confirm your registration, we would require you and one guest of your choice to provide us with:" "ඔබගේ ලියාපදිංචිය තහවුරු කිරීම සඳහා, පහත සඳහන් දෑ අපට ලබා දෙන ලෙස අපි ඔබෙන් සහ ඔබ කැමති එක් අමුත්තෙකුගෙන් ඉල්ලා සිටිමු:"

# 11. "A scanned copy of your international passports" (paragraph 20)
ReplaceInParagraph 20 "A scanned copy of your international passports" "ඔබගේ ජාත්‍යන්තර ගමන් බලපත්‍රවල පරිලෝකනය කළ පිටපතක්"

# 12. "Covid-19 vaccination certificates" (paragraph 21)
ReplaceInParagraph 21 "Covid-19 vaccination certificates" "Covid-19 එන්නත් සහතික"

# 13. "Send my details" (paragraph 23)
ReplaceInParagraph 23 "Send my details" "මගේ විස්තර යැවන්න"

# 14. "Your country manager will be in touch..." (paragraph 27)
ReplaceInParagraph 27 "Your country manager will be in touch to confirm your booking or request any other relevant details. " "ඔබේ වෙන් කිරීම තහවුරු කිරීමට හෝ වෙනත් අදාළ විස්තර ඉල්ලා සිටීමට ඔබේ රටේ කළමනාකරු ඔබ සමඟ​ සම්බන්ධ වනු ඇත. "

# 15. "Our event package offers you and your guest: " (paragraph 28)
ReplaceInParagraph 28 "Our event package offers you and your guest: " "අපගේ සිදුවීම් පැකේජයෙන් ඔබට සහ ඔබේ අමුත්තාට​ පහත සඳහන් දේ පිරිනමනු ඇත: "

# 16. "Travel insurance " (paragraph 30)
ReplaceInParagraph 30 "Travel insurance " "ගමන් රක්ෂණය "

# 17. "Airport – Hotel – Airport transfer " (paragraph 31)
ReplaceInParagraph 31 "Airport – Hotel – Airport transfer " "ගුවන් තොටුපළ - හෝටලය - ගුවන් තොටුපළ මාරු කිරීම "

# 18. "One hotel room for you and your guest / Two hotel rooms for you and your guest" (paragraph 32)
ReplaceInParagraph 32 "One hotel room for you and your guest / Two hotel rooms for you and your guest" "ඔබට සහ ඔබේ අමුත්තන්ට එක් හෝටල් කාමරයක් / ඔබට සහ ඔබේ අමුත්තන්ට හෝටල් කාමර දෙකක්"

# 19. "Check-in" + " on " (paragraph 33)
ReplaceInParagraph 33 "Check-in" "ඇතුළු වන්න"
ReplaceInParagraph 33 " on " " "

# 20. "Check-out" + " on " (paragraph 34)
ReplaceInParagraph 34 "Check-out" "පිටවීම"
ReplaceInParagraph 34 " on " " "

# 21. "Meals (Breakfast, lunch, and dinner)" (paragraph 35)
ReplaceInParagraph 35 "Meals (Breakfast, lunch, and dinner)" "ආහාර (උදෑසන, දිවා ආහාරය සහ රාත්‍රී ආහාරය)"

# 22-25: paragraph 40, several runs
ReplaceInParagraph 40 "If you have any questions, please contact your country manager, " "ඔබට කිසියම් ප්‍රශ්නයක් ඇත්නම්, කරුණාකර ඔබගේ රටේ කළමනාකරු "
ReplaceInParagraph 40 ", at " " ට "
ReplaceInParagraph 40 " or " " හෝ "
ReplaceInParagraph 40 " (WhatsApp). " " (WhatsApp) මගින් අමතන්න. "

# 26. "We look forward to seeing you soon." (paragraph 42)
ReplaceInParagraph 42 "We look forward to seeing you soon." "අපි ඉක්මනින් ඔබව දැකීමට බලාපොරොත්තු වෙමු."

Write-Output "all done"
